$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asistencia")
$ws.Range("G1").Value = 45793
